$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 926.2222
$ws.Range("J28").Value = 1099.8
$ws.Range("L28").Value = 1099.8
$ws.Range("N28").Value = -2069.8

$ws.Range("H29").Value = 650
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 3000
$ws.Range("N29").Value = -3562

$ws.Range("H31").Value = 3984.3333
$ws.Range("I31").Value = 3984.3333
$ws.Range("K31").Value = 11952.9999
$ws.Range("M31").Value = -11722.9999

$ws.Range("H38").Value = 2090.5557
$ws.Range("J38").Value = 676.75
$ws.Range("L38").Value = 2030.25
$ws.Range("N38").Value = -2774.25

$ws.Range("H40").Value = 990
$ws.Range("I40").Value = 990
$ws.Range("K40").Value = 990
$ws.Range("M40").Value = -815

$ws.Range("H58").Value = 716.25

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null

$ws.Range("H98").Value = 2433.3333
$ws.Range("I98").Value = 2433.3333
$ws.Range("K98").Value = 2433.3333
$ws.Range("M98").Value = -935.3332999999998

$ws.Range("H122").Value = 2433.3333
$ws.Range("I122").Value = 2433.3333
$ws.Range("K122").Value = 7299.999899999999
$ws.Range("M122").Value = -4849.999899999999

$ws.Range("H132").Value = 2739
$ws.Range("I132").Value = 2739
$ws.Range("K132").Value = 8217
$ws.Range("M132").Value = -5687

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2999.5
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1623

$ws.Range("H74").Value = 1032.8889
$ws.Range("I74").Value = 1074.5
$ws.Range("J74").Value = 700
$ws.Range("K74").Value = 1074.5
$ws.Range("L74").Value = 700
$ws.Range("M74").Value = -200.5
$ws.Range("N74").Value = -2448

$ws.Range("H77").Value = 1032.8889
$ws.Range("I77").Value = 1074.5
$ws.Range("J77").Value = 700
$ws.Range("K77").Value = 5372.5
$ws.Range("L77").Value = 3500
$ws.Range("M77").Value = -1004.5
$ws.Range("N77").Value = -12236

$ws.Range("H88").Value = 360.25
$ws.Range("I88").Value = 268.85715
$ws.Range("J88").Value = 1000
$ws.Range("K88").Value = 268.85715
$ws.Range("L88").Value = 1000
$ws.Range("M88").Value = 137.14285
$ws.Range("N88").Value = -1812

$ws.Range("H91").Value = 360.25
$ws.Range("I91").Value = 268.85715
$ws.Range("J91").Value = 1000
$ws.Range("K91").Value = 268.85715
$ws.Range("L91").Value = 1000
$ws.Range("M91").Value = 1135.14285
$ws.Range("N91").Value = -3808

$ws.Range("H132").Value = 2326.8667
$ws.Range("I132").Value = 2736.6365
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 8209.9095
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -5679.9095
$ws.Range("N132").Value = -8660

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 72541.5
$ws.Range("I54").Value = 72541.5
$ws.Range("K54").Value = 72541.5
$ws.Range("M54").Value = -72057.5

$ws.Range("H86").Value = 2503.5
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 2007
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 2007
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -4253

$ws.Range("H89").Value = 2503.5
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 2007
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 10035
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -21267

$ws.Range("H107").Value = 1663.125
$ws.Range("I107").Value = 1507.8572
$ws.Range("K107").Value = 1507.8572
$ws.Range("M107").Value = 412.1428000000001

$ws.Range("H134").Value = 1521.2222
$ws.Range("I134").Value = 1336.375
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 4009.125
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -1474.125
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 413.83334
$ws.Range("I107").Value = 419.6
$ws.Range("J107").Value = 385
$ws.Range("K107").Value = 419.6
$ws.Range("L107").Value = 385
$ws.Range("M107").Value = 1500.4
$ws.Range("N107").Value = -4225

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -188

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = $null

$ws.Range("H23").Value = 258.77777
$ws.Range("J23").Value = 431.8
$ws.Range("L23").Value = 1295.4
$ws.Range("N23").Value = -1765.4

$ws.Range("H29").Value = 25
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 25
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 75
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = -629

$ws.Range("H33").Value = 120.4
$ws.Range("I33").Value = 121.666664
$ws.Range("J33").Value = 118.5
$ws.Range("K33").Value = 729.999984
$ws.Range("L33").Value = 711
$ws.Range("M33").Value = -446.999984
$ws.Range("N33").Value = -1277

$ws.Range("H46").Value = 1650
$ws.Range("I46").Value = 950
$ws.Range("K46").Value = 2850
$ws.Range("M46").Value = -2759

$ws.Range("H55").Value = 2501.25
$ws.Range("I55").Value = 2750
$ws.Range("J55").Value = 2252.5
$ws.Range("K55").Value = 8250
$ws.Range("L55").Value = 6757.5
$ws.Range("M55").Value = -8073
$ws.Range("N55").Value = -7111.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 44999.5
$ws.Range("J39").Value = 44999.5
$ws.Range("L39").Value = 44999.5
$ws.Range("N39").Value = -46063.5

$ws.Range("H123").Value = 90098
$ws.Range("I123").Value = 70296
$ws.Range("J123").Value = 99999
$ws.Range("K123").Value = 70296
$ws.Range("L123").Value = 99999
$ws.Range("M123").Value = -67846
$ws.Range("N123").Value = -104899

$ws.Range("H132").Value = 3832
$ws.Range("I132").Value = 3832
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11496
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8966
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 18001500
$ws.Range("J42").Value = 35000000
$ws.Range("L42").Value = 35000000
$ws.Range("N42").Value = -35001126

$ws.Range("H49").Value = 18001500
$ws.Range("J49").Value = 35000000
$ws.Range("L49").Value = 35000000
$ws.Range("N49").Value = -35000294

$ws.Range("H137").Value = 47700
$ws.Range("I137").Value = 47700
$ws.Range("K137").Value = 47700
$ws.Range("M137").Value = -42600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = $null

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = $null

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -8060
